$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# This sheet is a "Saldo" (balance) export with columns: Conta | Nome | Saldo
# starting at row 1 (header) and row 2 onward holding data, terminated by a
# blank row followed by a "Filtros aplicados" note.
#
# The edit:
#   * Removes the account 005685089 / CARNEIRO / 15000 row entirely.
#   * Removes four rows: 004752461/SERGIO/6000, 001000882/AYRTON/3000,
#     005554830/PAULO/2867.29 (old entry), 002064834/RAFAELA/2000,
#     004487140/VALMIR/1173.96.
#   * Adds a new row for account 005554830 / PAULO with an updated balance
#     of 23867.29, positioned immediately before the 004364200/BLOCO row
#     (i.e. Paulo's balance entry moved up and was corrected upward).
# ---------------------------------------------------------------------------

function Get-LastDataRow {
    # Row 1 is the header ("Conta"/"Nome"/"Saldo"); data runs until the
    # first completely blank "Conta" cell (the sheet ends with a blank row
    # followed by a footnote row).
    $r = 2
    while ($ws.Cells.Item($r, 1).Text -ne "") {
        $r = $r + 1
    }
    return $r - 1
}

function Find-RowByAccount([string]$account, [int]$startRow, [int]$endRow) {
    for ($r = $startRow; $r -le $endRow; $r++) {
        if ($ws.Cells.Item($r, 1).Text -eq $account) {
            return $r
        }
    }
    return -1
}

$lastRow = Get-LastDataRow

# --- Remove the rows that disappear from the export -----------------------
$accountsToDelete = @("005685089", "004752461", "001000882", "005554830", "002064834", "004487140")

foreach ($acct in $accountsToDelete) {
    $lastRow = Get-LastDataRow
    $row = Find-RowByAccount $acct 2 $lastRow
    if ($row -ge 2) {
        $ws.Rows.Item($row).Delete()
    }
}

# --- Insert the new / corrected Paulo row right before the BLOCO row ------
$lastRow = Get-LastDataRow
$blocoRow = Find-RowByAccount "004364200" 2 $lastRow

$ws.Rows.Item($blocoRow).Insert()
$ws.Cells.Item($blocoRow, 1).Value = "'005554830"
$ws.Cells.Item($blocoRow, 2).Value = "PAULO"
$ws.Cells.Item($blocoRow, 3).Value = 23867.29
